$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17
$ws.Range("C2").Value = 6
$ws.Range("E2").Value = 6387103.287225001
$ws.Range("G2").Value = 45183
$ws.Range("H2").Value = 28
$ws.Range("I2").Value = 13
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 0.464
$ws.Range("L2").Value = 0.536
$ws.Range("M2").Value = -0.07000000000000001
$ws.Range("N2").Value = 428826.258
$ws.Range("O2").Value = 0.188
$ws.Range("P2").Value = -92574.789
$ws.Range("Q2").Value = -0.052
$ws.Range("R2").Value = -236304.743
$ws.Range("S2").Value = 2330591.91
$ws.Range("T2").Value = 2.33059191015625
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 6
$ws.Range("W2").Value = 53.2962962962963
